$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 75, shifting rows 75-109 down to 76-110
$ws.Rows.Item(75).Insert()

# Populate the new row 75 with the weekly entry (copy constant columns from row 76 pattern)
$ws.Cells.Item(75, 1).Value = 7
$ws.Cells.Item(75, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(75, 3).Value = "Ñuble"
$ws.Cells.Item(75, 4).Value = 45134
$ws.Cells.Item(75, 4).NumberFormat = $ws.Cells.Item(76, 4).NumberFormat
$ws.Cells.Item(75, 5).Value = 16
$ws.Cells.Item(75, 6).Value = 100112013
$ws.Cells.Item(75, 7).Value = "Alcachofa"
$ws.Cells.Item(75, 8).Value = "Madrigal"
$ws.Cells.Item(75, 9).Value = "Primera"
$ws.Cells.Item(75, 10).Value = 60
$ws.Cells.Item(75, 11).Value = 15000
$ws.Cells.Item(75, 12).Value = 15000
$ws.Cells.Item(75, 13).Value = 15000
$ws.Cells.Item(75, 14).Value = "`$/caja 40 unidades"
$ws.Cells.Item(75, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(75, 16).Value = 375
$ws.Cells.Item(75, 17).Value = 40
$ws.Cells.Item(75, 18).Value = "Hortaliza"
